# Update roc curve formatting
# - Add an "AUC CI" column to the results table.
# - Replace the data rows with the refreshed metrics (new algo/task ordering,
#   new NeuralNetBinaryClassifier rows, updated numbers, and AUC CI strings).
# - Re-flow every column width from 1440 dxa to 1234 dxa.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Final table contents (header row + 9 data rows), 7 columns wide.
$data = @(
    @('Algo Name','Task','Accuracy','Specificity','Sensitivity','AUC','AUC CI'),
    @('DecisionTreeClassifier','LOS','0.84','0.88','0.55','0.81','0.7876978709831949-0.8389850498964259'),
    @('RandomForestClassifier','LOS','0.84','0.90','0.43','0.79','0.7766429755234755-0.8090514570550523'),
    @('DecisionTreeClassifier','DIED','0.94','0.94','0.68','0.89','0.8365671357836144-0.9463630398262148'),
    @('NeuralNetBinaryClassifier','anastomotic_leak','0.01','0.00','1.00','0.88','0.8425944328020263-0.9222740158910612'),
    @('RandomForestClassifier','anastomotic_leak','0.98','0.99','0.03','0.83','0.7760623792933818-0.890945327864394'),
    @('NeuralNetBinaryClassifier','LOS','0.13','0.00','1.00','0.84','0.8211942636312464-0.8500799841397157'),
    @('DecisionTreeClassifier','anastomotic_leak','0.93','0.94','0.39','0.84','0.7620435291061813-0.9136384845908263'),
    @('NeuralNetBinaryClassifier','DIED','0.01','0.00','1.00','0.90','0.846252187840221-0.9571426780726114'),
    @('RandomForestClassifier','DIED','0.98','0.99','0.16','0.90','0.8563649716654991-0.9437602318194811')
)

# Add the 7th ("AUC CI") column - Word appends it after the last column.
if ($t.Columns.Count -lt 7) {
    $t.Columns.Add() | Out-Null
}

# Write every cell's text to match the new table contents.
for ($r = 1; $r -le $data.Count; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Count; $c++) {
        $t.Cell($r, $c).Range.Text = $row[$c - 1]
    }
}

# Re-flow column widths: 1440 dxa (72pt) -> 1234 dxa (61.7pt).
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $t.Columns.Item($c).Width = 61.7
}
